# The workbook tracks one price record per row (Berenjena / Mercado
# Mayorista Lo Valledor de Santiago). A new weekly record needs to be
# inserted as row 87, pushing the existing row 87 (and everything below
# it) down by one row, all the way to the end of the sheet (old row 165
# becomes row 166).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 87; Excel shifts rows 87:165 down to
# 88:166 and grows the used range automatically.
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new record's data.
$ws.Cells.Item(87, 1).Value2 = 6
$ws.Cells.Item(87, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(87, 3).Value2 = "Metropolitana"
$ws.Cells.Item(87, 4).Value2 = 44586
$ws.Cells.Item(87, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(87, 5).Value2 = 13
$ws.Cells.Item(87, 6).Value2 = 100112001
$ws.Cells.Item(87, 7).Value2 = "Berenjena"
$ws.Cells.Item(87, 8).Value2 = "Sin especificar"
$ws.Cells.Item(87, 9).Value2 = "Primera"
$ws.Cells.Item(87, 10).Value2 = 230
$ws.Cells.Item(87, 11).Value2 = 10000
$ws.Cells.Item(87, 12).Value2 = 12000
$ws.Cells.Item(87, 13).Value2 = 10870
$ws.Cells.Item(87, 14).Value2 = "$/caja 50 unidades"
$ws.Cells.Item(87, 15).Value2 = "Provincia de Huasco"
$ws.Cells.Item(87, 16).Value2 = 217
$ws.Cells.Item(87, 17).Value2 = 50
$ws.Cells.Item(87, 18).Value2 = "Hortaliza"
